$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook was re-saved by Excel (on a machine/session with slightly
# different display/font metrics), which nudges the auto-computed column A
# "best fit" width and the wrap-text auto row heights for rows 2-4 by a
# small amount. Reproduce the resulting geometry as closely as this object
# model allows.

$ws.Columns.Item(1).ColumnWidth = 13.75

$ws.Rows.Item(2).RowHeight = 173.25
$ws.Rows.Item(3).RowHeight = 173.25
$ws.Rows.Item(4).RowHeight = 220.5
